$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 400
$ws.Range("I2").Value = 400
$ws.Range("K2").Value = 400
$ws.Range("M2").Value = -287
$ws.Range("H11").Value = 96.88
$ws.Range("I11").Value = 96.88
$ws.Range("K11").Value = 96.88
$ws.Range("M11").Value = 43.12
$ws.Range("H46").Value = 2395
$ws.Range("I46").Value = 100
$ws.Range("J46").Value = 3160
$ws.Range("K46").Value = 300
$ws.Range("L46").Value = 9480
$ws.Range("M46").Value = -181
$ws.Range("N46").Value = -9718
$ws.Range("H60").Value = 2395
$ws.Range("I60").Value = 100
$ws.Range("J60").Value = 3160
$ws.Range("K60").Value = 300
$ws.Range("L60").Value = 9480
$ws.Range("M60").Value = 184
$ws.Range("N60").Value = -10448
$ws.Range("H69").Value = 3328.3333
$ws.Range("I69").Value = 3323.1667
$ws.Range("K69").Value = 9969.500100000001
$ws.Range("M69").Value = -9095.500100000001
$ws.Range("H72").Value = 3328.3333
$ws.Range("I72").Value = 3323.1667
$ws.Range("K72").Value = 29908.5003
$ws.Range("M72").Value = -25540.5003
$ws.Range("H74").Value = 2976.7058
$ws.Range("I74").Value = 2571.1667
$ws.Range("J74").Value = 3950
$ws.Range("K74").Value = 2571.1667
$ws.Range("L74").Value = 3950
$ws.Range("M74").Value = -1635.1667
$ws.Range("N74").Value = -5822
$ws.Range("H77").Value = 2976.7058
$ws.Range("I77").Value = 2571.1667
$ws.Range("J77").Value = 3950
$ws.Range("K77").Value = 12855.8335
$ws.Range("L77").Value = 19750
$ws.Range("M77").Value = -8175.833500000001
$ws.Range("N77").Value = -29110
$ws.Range("H94").Value = 3968.5
$ws.Range("I94").Value = 921.25
$ws.Range("K94").Value = 921.25
$ws.Range("M94").Value = -470.25

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21743.44
$ws.Range("I32").Value = 22488
$ws.Range("J32").Value = 14000
$ws.Range("K32").Value = 22488
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = -22201
$ws.Range("N32").Value = -14574
$ws.Range("H61").Value = 1747.3206
$ws.Range("I61").Value = 1748.9524
$ws.Range("J61").Value = 1740.4667
$ws.Range("K61").Value = 1748.9524
$ws.Range("L61").Value = 1740.4667
$ws.Range("M61").Value = -1536.9524
$ws.Range("N61").Value = -2164.4667
$ws.Range("H132").Value = 1117.5077
$ws.Range("I132").Value = 914.68256
$ws.Range("K132").Value = 2744.04768
$ws.Range("M132").Value = -214.0476799999997
$ws.Range("H136").Value = 1747.3206
$ws.Range("I136").Value = 1748.9524
$ws.Range("J136").Value = 1740.4667
$ws.Range("K136").Value = 5246.857199999999
$ws.Range("L136").Value = 5221.4001
$ws.Range("M136").Value = -2696.857199999999
$ws.Range("N136").Value = -10321.4001
$ws.Range("H138").Value = 53357.715
$ws.Range("J138").Value = 53357.715
$ws.Range("L138").Value = 53357.715
$ws.Range("N138").Value = -63637.715

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2119.724
$ws.Range("I86").Value = 1969.1818
$ws.Range("J86").Value = 2592.8572
$ws.Range("K86").Value = 1969.1818
$ws.Range("L86").Value = 2592.8572
$ws.Range("M86").Value = -846.1818000000001
$ws.Range("N86").Value = -4838.8572
$ws.Range("H89").Value = 2119.724
$ws.Range("I89").Value = 1969.1818
$ws.Range("J89").Value = 2592.8572
$ws.Range("K89").Value = 9845.909
$ws.Range("L89").Value = 12964.286
$ws.Range("M89").Value = -4229.909
$ws.Range("N89").Value = -24196.286
$ws.Range("H118").Value = 7877.5
$ws.Range("J118").Value = 7877.5
$ws.Range("L118").Value = 7877.5
$ws.Range("N118").Value = -11191.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6804.7417
$ws.Range("I31").Value = 2257.9565
$ws.Range("J31").Value = 19876.75
$ws.Range("K31").Value = 2257.9565
$ws.Range("L31").Value = 19876.75
$ws.Range("M31").Value = -1962.9565
$ws.Range("N31").Value = -20466.75
$ws.Range("H34").Value = 6804.7417
$ws.Range("I34").Value = 2257.9565
$ws.Range("J34").Value = 19876.75
$ws.Range("K34").Value = 2257.9565
$ws.Range("L34").Value = 19876.75
$ws.Range("M34").Value = -2055.9565
$ws.Range("N34").Value = -20280.75
$ws.Range("H132").Value = 3417.0386
$ws.Range("I132").Value = 1230.6818
$ws.Range("J132").Value = 15442
$ws.Range("K132").Value = 3692.0454
$ws.Range("L132").Value = 46326
$ws.Range("M132").Value = -1162.0454
$ws.Range("N132").Value = -51386

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1557
$ws.Range("I68").Value = 2500
$ws.Range("J68").Value = 1179.8
$ws.Range("K68").Value = 7500
$ws.Range("L68").Value = 3539.4
$ws.Range("M68").Value = -6689
$ws.Range("N68").Value = -5161.4
$ws.Range("H71").Value = 1557
$ws.Range("I71").Value = 2500
$ws.Range("J71").Value = 1179.8
$ws.Range("K71").Value = 22500
$ws.Range("L71").Value = 10618.2
$ws.Range("M71").Value = -18444
$ws.Range("N71").Value = -18730.2
$ws.Range("H122").Value = 1264.3636
$ws.Range("I122").Value = 1347.091
$ws.Range("J122").Value = 1181.6364
$ws.Range("K122").Value = 12123.819
$ws.Range("L122").Value = 10634.7276
$ws.Range("M122").Value = -9673.819
$ws.Range("N122").Value = -15534.7276
$ws.Range("H132").Value = 1450.5927
$ws.Range("I132").Value = 721.3333
$ws.Range("J132").Value = 1658.9524
$ws.Range("K132").Value = 6491.9997
$ws.Range("L132").Value = 14930.5716
$ws.Range("M132").Value = -3961.9997
$ws.Range("N132").Value = -19990.5716

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2837.4375
$ws.Range("I80").Value = 2833.3333
$ws.Range("J80").Value = 2838.3845
$ws.Range("K80").Value = 2833.3333
$ws.Range("L80").Value = 2838.3845
$ws.Range("M80").Value = -1835.3333
$ws.Range("N80").Value = -4834.3845
$ws.Range("H83").Value = 2837.4375
$ws.Range("I83").Value = 2833.3333
$ws.Range("J83").Value = 2838.3845
$ws.Range("K83").Value = 14166.6665
$ws.Range("L83").Value = 14191.9225
$ws.Range("M83").Value = -9174.666499999999
$ws.Range("N83").Value = -24175.9225
$ws.Range("H133").Value = 59799.5
$ws.Range("J133").Value = 59799.5
$ws.Range("L133").Value = 59799.5
$ws.Range("N133").Value = -69919.5
$ws.Range("H141").Value = 66652.86
$ws.Range("J141").Value = 66652.86
$ws.Range("L141").Value = 66652.86
$ws.Range("N141").Value = -77012.86

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2497.1155
$ws.Range("I68").Value = 2182.2104
$ws.Range("K68").Value = 2182.2104
$ws.Range("M68").Value = -1433.2104
$ws.Range("H71").Value = 2497.1155
$ws.Range("I71").Value = 2182.2104
$ws.Range("K71").Value = 10911.052
$ws.Range("M71").Value = -7167.052
$ws.Range("H82").Value = 1734.7
$ws.Range("I82").Value = 1705.2222
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 1705.2222
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -1344.2222
$ws.Range("N82").Value = -2722
$ws.Range("H85").Value = 1734.7
$ws.Range("I85").Value = 1705.2222
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 1705.2222
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -457.2221999999999
$ws.Range("N85").Value = -4496

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1922.2449
$ws.Range("I126").Value = 1939.7727
$ws.Range("K126").Value = 5819.3181
$ws.Range("L126").Value = 5723.889
$ws.Range("M126").Value = -3349.3181
$ws.Range("N126").Value = -10663.889
